# GP05MOAS-GL276 calibration sheet update
# redmine #9229 - Calibration sheets added/changed for GP05MOAS gliders
# GL276, GL361-GL365, GL453, GL523, GL525, GL537, PG514, PG515.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Moorings": update the mooring deployment record ----
$moorings = $wb.Worksheets.Item("Moorings")

# Anchor Launch Date: 13-Jun-2014 -> 15-Jun-2014 (serial date, preserves existing date format)
$moorings.Range("D2").Value = 41805
# Anchor Launch Time: 0:00 -> 5:30 AM (serial fraction of a day, preserves existing time format)
$moorings.Range("E2").Value = 5.5/24
# Recover Date: (blank) -> 02-Jun-2015
$moorings.Range("F2").Value = 42157
# Mooring Serial Number: MV1404 -> MV-1404
$moorings.Range("J2").Value = "MV-1404"

# Flag the edited cells with a blue font color (matches author convention)
$moorings.Range("D2:F2").Font.Color = 16711680  # RGB(0,0,255) blue
$moorings.Range("J2").Font.Color = 16711680      # RGB(0,0,255) blue

# ---- Sheet "Asset_Cal_Info": update a calibration coefficient value ----
$assetCal = $wb.Worksheets.Item("Asset_Cal_Info")

# CC_angular_resolution value: 1.13 -> 1.096
$assetCal.Range("F6").Value = 1.096

# Flag the edited cell with blue font + yellow highlight
$assetCal.Range("F6").Font.Color = 16711680  # RGB(0,0,255) blue
$assetCal.Range("F6").Interior.Color = 65535  # RGB(255,255,0) yellow
